$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($ws, $ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-CellText $ws "D2" "31.168.00"
Set-CellText $ws "E2" "  +1.44%  "

# Row 3
Set-CellText $ws "D3" "1.961.66"
Set-CellText $ws "E3" "  +0.76%  "

# Row 4
Set-CellText $ws "D4" "1.003"
Set-CellText $ws "E4" "  +0.25%  "

# Row 5
Set-CellText $ws "D5" "246.30"
Set-CellText $ws "E5" "  -0.50%  "

# Row 6
Set-CellText $ws "E6" "  +0.23%  "

# Row 7
Set-CellText $ws "D7" "0.4915"
Set-CellText $ws "E7" "  +1.55%  "

# Row 8
Set-CellText $ws "D8" "0.3008"
Set-CellText $ws "E8" "  +1.86%  "

# Row 9
Set-CellText $ws "B9" "Dogecoin"
Set-CellText $ws "C9" "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-CellText $ws "D9" "0.06899"
Set-CellText $ws "E9" "  +1.16%  "

# Row 10
Set-CellText $ws "B10" "Litecoin"
Set-CellText $ws "C10" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-CellText $ws "D10" "109.01"
Set-CellText $ws "E10" "  -3.15%  "

# Row 11
Set-CellText $ws "D11" "19.24"
Set-CellText $ws "E11" "  -1.38%  "

# Row 12
Set-CellText $ws "B12" "TRON"
Set-CellText $ws "C12" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-CellText $ws "D12" "0.07785"
Set-CellText $ws "E12" "  +1.72%  "

# Row 13
Set-CellText $ws "B13" "WrappedEther"
Set-CellText $ws "C13" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-CellText $ws "D13" "1.931.18"
Set-CellText $ws "E13" "  -0.83%  "

# Row 14
Set-CellText $ws "B14" "Polkadot"
Set-CellText $ws "C14" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-CellText $ws "D14" "5.475"
Set-CellText $ws "E14" "  -1.29%  "

# Row 15
Set-CellText $ws "B15" "Polygon"
Set-CellText $ws "C15" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-CellText $ws "D15" "0.7136"
Set-CellText $ws "E15" "  +3.16%  "

# Row 16
Set-CellText $ws "B16" "BitcoinCash"
Set-CellText $ws "C16" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-CellText $ws "D16" "284.93"
Set-CellText $ws "E16" "  -3.96%  "

# Row 17
Set-CellText $ws "B17" "WrappedBTC"
Set-CellText $ws "C17" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-CellText $ws "D17" "31.047.17"
Set-CellText $ws "E17" "  +1.05%  "

# Row 18
Set-CellText $ws "B18" "Avalanche"
Set-CellText $ws "C18" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-CellText $ws "D18" "13.31"
Set-CellText $ws "E18" "  -0.45%  "

# Row 19
Set-CellText $ws "B19" "ShibaInu"
Set-CellText $ws "C19" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-CellText $ws "D19" "0.000007784"
Set-CellText $ws "E19" "  +1.04%  "

# Row 20
Set-CellText $ws "B20" "Dai"
Set-CellText $ws "C20" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-CellText $ws "D20" "1.003"
Set-CellText $ws "E20" "  +0.24%  "

# Row 21
Set-CellText $ws "B21" "WrappedliquidstakedEther2.0"
Set-CellText $ws "C21" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-CellText $ws "D21" "2.177.84"
Set-CellText $ws "E21" "  -0.97%  "

# Row 22
Set-CellText $ws "B22" "Uniswap"
Set-CellText $ws "C22" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-CellText $ws "D22" "5.518"
Set-CellText $ws "E22" "  -3.11%  "

# Row 23
Set-CellText $ws "B23" "BinanceUSD"
Set-CellText $ws "C23" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-CellText $ws "D23" "1.003"
Set-CellText $ws "E23" "  +0.26%  "

# Row 24
Set-CellText $ws "B24" "Chainlink"
Set-CellText $ws "C24" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-CellText $ws "D24" "6.563"
Set-CellText $ws "E24" "  +0.09%  "

# Row 25
Set-CellText $ws "B25" "Cosmos"
Set-CellText $ws "C25" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-CellText $ws "D25" "9.840"
Set-CellText $ws "E25" "  +0.55%  "

# Row 26
Set-CellText $ws "B26" "Monero"
Set-CellText $ws "C26" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-CellText $ws "D26" "169.93"
Set-CellText $ws "E26" "  +1.10%  "

# Row 27
Set-CellText $ws "B27" "EthereumClassic"
Set-CellText $ws "C27" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-CellText $ws "D27" "20.19"
Set-CellText $ws "E27" "  -0.79%  "

# Row 28
Set-CellText $ws "B28" "LidoDAOToken"
Set-CellText $ws "C28" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-CellText $ws "D28" "2.233"
Set-CellText $ws "E28" "  +2.44%  "

# Row 29
Set-CellText $ws "B29" "Stellar"
Set-CellText $ws "C29" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-CellText $ws "D29" "0.1051"
Set-CellText $ws "E29" "  -3.65%  "

# Row 30
Set-CellText $ws "B30" "Toncoin"
Set-CellText $ws "C30" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-CellText $ws "D30" "1.438"
Set-CellText $ws "E30" "  +0.06%  "

# Row 31
Set-CellText $ws "B31" "PancakeSwap"
Set-CellText $ws "C31" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-CellText $ws "D31" "1.586"

# Row 32
Set-CellText $ws "B32" "Filecoin"
Set-CellText $ws "C32" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-CellText $ws "D32" "4.624"
Set-CellText $ws "E32" "  -2.79%  "

# Row 33
Set-CellText $ws "B33" "InternetComputer(DFINITY)"
Set-CellText $ws "C33" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-CellText $ws "D33" "4.487"
Set-CellText $ws "E33" "  +1.44%  "

# Row 34
Set-CellText $ws "B34" "Hedera"
Set-CellText $ws "C34" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-CellText $ws "D34" "0.04995"
Set-CellText $ws "E34" "  -1.54%  "

# Row 35
Set-CellText $ws "B35" "ImmutableX"
Set-CellText $ws "C35" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-CellText $ws "D35" "0.7632"
Set-CellText $ws "E35" "  -2.07%  "

# Row 36
Set-CellText $ws "B36" "ARBITRUM"
Set-CellText $ws "C36" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-CellText $ws "D36" "1.187"
Set-CellText $ws "E36" "  +2.16%  "

# Row 37
Set-CellText $ws "B37" "HuobiToken"
Set-CellText $ws "C37" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-CellText $ws "D37" "2.739"
Set-CellText $ws "E37" "  +0.11%  "

# Row 38
Set-CellText $ws "D38" "0.02050"
Set-CellText $ws "E38" "  -1.18%  "

# Row 39
Set-CellText $ws "B39" "MXToken"
Set-CellText $ws "C39" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-CellText $ws "D39" "2.712"
Set-CellText $ws "E39" "  +0.43%  "

# Row 40
Set-CellText $ws "B40" "RenderToken"
Set-CellText $ws "C40" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-CellText $ws "D40" "2.180"
Set-CellText $ws "E40" "  +6.55%  "

# Row 41
Set-CellText $ws "B41" "FraxShare"
Set-CellText $ws "C41" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-CellText $ws "D41" "6.486"
Set-CellText $ws "E41" "  +8.88%  "

# Row 42
Set-CellText $ws "B42" "TheSandbox"
Set-CellText $ws "C42" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-CellText $ws "D42" "0.4535"
Set-CellText $ws "E42" "  +1.76%  "

# Row 43
Set-CellText $ws "B43" "Aave"
Set-CellText $ws "C43" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-CellText $ws "D43" "73.27"
Set-CellText $ws "E43" "  +3.24%  "

# Row 44
Set-CellText $ws "B44" "TrustWalletToken"
Set-CellText $ws "C44" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-CellText $ws "D44" "0.8865"
Set-CellText $ws "E44" "  +1.35%  "

# Row 45
Set-CellText $ws "D45" "109.58"
Set-CellText $ws "E45" "  -1.26%  "

# Row 46
Set-CellText $ws "D46" "8.177"
Set-CellText $ws "E46" "  +10.55%  "

# Row 47
Set-CellText $ws "B47" "PaxDollar"
Set-CellText $ws "C47" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-CellText $ws "D47" "1.003"
Set-CellText $ws "E47" "  -0.02%  "

# Row 48
Set-CellText $ws "B48" "EnergySwap"
Set-CellText $ws "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-CellText $ws "D48" "9.484"
Set-CellText $ws "E48" "  +0.24%  "

# Row 49
Set-CellText $ws "B49" "Maker"
Set-CellText $ws "C49" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-CellText $ws "D49" "965.39"
Set-CellText $ws "E49" "  +5.82%  "

# Row 50
Set-CellText $ws "D50" "0.1272"
Set-CellText $ws "E50" "  +1.70%  "

# Row 51
Set-CellText $ws "B51" "WOONetwork"
Set-CellText $ws "C51" "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
Set-CellText $ws "D51" "0.2607"
Set-CellText $ws "E51" "  +2.25%  "

